# Update the "Förändrad" (Changed) date column (C) for rows 2-11
# from serial 45207 (2023-10-08) to serial 45208 (2023-10-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
